# Auto-generated script to update FFXIV leve-profit market data cells
# Updates currentAveragePrice / LevePriceNQ / LevePriceHQ / LeveProfit columns
# across the ALC, ARM, CRP, CUL, GSM, LTW, WVR sheets to reflect refreshed
# market-board pricing pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 1098.75
$ws.Range("I46").Value = 965.3333
$ws.Range("K46").Value = 2895.9999
$ws.Range("M46").Value = -2776.9999

$ws.Range("H60").Value = 1098.75
$ws.Range("I60").Value = 965.3333
$ws.Range("K60").Value = 2895.9999
$ws.Range("M60").Value = -2411.9999

$ws.Range("H110").Value = 15583
$ws.Range("J110").Value = 15583
$ws.Range("L110").Value = 15583
$ws.Range("N110").Value = -23763

$ws.Range("H137").Value = 3200.516
$ws.Range("I137").Value = 2226.238
$ws.Range("J137").Value = 5246.5
$ws.Range("K137").Value = 6678.714
$ws.Range("L137").Value = 15739.5
$ws.Range("M137").Value = -4128.714
$ws.Range("N137").Value = -20839.5


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 835.5
$ws.Range("I32").Value = 677.08234
$ws.Range("K32").Value = 677.08234
$ws.Range("M32").Value = -390.08234

$ws.Range("H122").Value = 949467.9
$ws.Range("I122").Value = 2455.2222
$ws.Range("K122").Value = 7365.6666
$ws.Range("M122").Value = -4915.6666


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1642
$ws.Range("I16").Value = 1347.4736
$ws.Range("J16").Value = 2441.4285
$ws.Range("K16").Value = 1347.4736
$ws.Range("L16").Value = 2441.4285
$ws.Range("M16").Value = -1060.4736
$ws.Range("N16").Value = -3015.4285

$ws.Range("H58").Value = 2765.3125
$ws.Range("I58").Value = 2226.9167
$ws.Range("K58").Value = 2226.9167
$ws.Range("M58").Value = -2023.9167

$ws.Range("H109").Value = 51677.4
$ws.Range("J109").Value = 51677.4
$ws.Range("L109").Value = 51677.4
$ws.Range("N109").Value = -53757.4

$ws.Range("H113").Value = 1642
$ws.Range("I113").Value = 1347.4736
$ws.Range("J113").Value = 2441.4285
$ws.Range("K113").Value = 1347.4736
$ws.Range("L113").Value = 2441.4285
$ws.Range("M113").Value = 822.5264
$ws.Range("N113").Value = -6781.4285

$ws.Range("H132").Value = 3310.625
$ws.Range("I132").Value = 2580.8333
$ws.Range("K132").Value = 7742.499899999999
$ws.Range("M132").Value = -5212.499899999999

$ws.Range("H134").Value = 28936.324
$ws.Range("I134").Value = 35260.926
$ws.Range("J134").Value = 4541.4287
$ws.Range("K134").Value = 105782.778
$ws.Range("L134").Value = 13624.2861
$ws.Range("M134").Value = -103247.778
$ws.Range("N134").Value = -18694.2861

$ws.Range("H136").Value = 2765.3125
$ws.Range("I136").Value = 2226.9167
$ws.Range("K136").Value = 6680.750100000001
$ws.Range("M136").Value = -4130.750100000001


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

$ws.Range("H28").Value = 5000
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 5000
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 15000
$ws.Range("M28").ClearContents()
$ws.Range("N28").Value = -15464

$ws.Range("H41").Value = 2526552.2
$ws.Range("I41").Value = 5051005
$ws.Range("J41").Value = 2099.5
$ws.Range("K41").Value = 15153015
$ws.Range("L41").Value = 6298.5
$ws.Range("M41").Value = -15152677
$ws.Range("N41").Value = -6974.5

$ws.Range("H43").Value = 1231.4
$ws.Range("I43").Value = 582.5
$ws.Range("J43").Value = 1664
$ws.Range("K43").Value = 1747.5
$ws.Range("L43").Value = 4992
$ws.Range("M43").Value = -1633.5
$ws.Range("N43").Value = -5220

$ws.Range("H97").Value = 256.6
$ws.Range("I97").Value = 256.6
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 769.8000000000001
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -273.8000000000001
$ws.Range("N97").ClearContents()

$ws.Range("H98").Value = 251.75
$ws.Range("I98").Value = 252.33333
$ws.Range("J98").Value = 250
$ws.Range("K98").Value = 756.99999
$ws.Range("L98").Value = 750
$ws.Range("M98").Value = 741.00001
$ws.Range("N98").Value = -3746

$ws.Range("H137").Value = 2559.6428
$ws.Range("I137").Value = 1976.9
$ws.Range("J137").Value = 4016.5
$ws.Range("K137").Value = 5930.700000000001
$ws.Range("L137").Value = 12049.5
$ws.Range("M137").Value = -830.7000000000007
$ws.Range("N137").Value = -22249.5


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1628972.6
$ws.Range("I80").Value = 4065572.2
$ws.Range("J80").Value = 4572.8887
$ws.Range("K80").Value = 4065572.2
$ws.Range("L80").Value = 4572.8887
$ws.Range("M80").Value = -4064574.2
$ws.Range("N80").Value = -6568.8887

$ws.Range("H83").Value = 1628972.6
$ws.Range("I83").Value = 4065572.2
$ws.Range("J83").Value = 4572.8887
$ws.Range("K83").Value = 20327861
$ws.Range("L83").Value = 22864.4435
$ws.Range("M83").Value = -20322869
$ws.Range("N83").Value = -32848.4435

$ws.Range("H106").Value = 29000
$ws.Range("J106").Value = 29000
$ws.Range("L106").Value = 29000
$ws.Range("N106").Value = -31524

$ws.Range("H113").Value = 27778916
$ws.Range("I113").Value = 41667124
$ws.Range("J113").Value = 2499.5
$ws.Range("K113").Value = 41667124
$ws.Range("L113").Value = 2499.5
$ws.Range("M113").Value = -41664954
$ws.Range("N113").Value = -6839.5

$ws.Range("H122").Value = 3552.4412
$ws.Range("I122").Value = 2492.6538
$ws.Range("K122").Value = 7477.9614
$ws.Range("M122").Value = -5027.9614

$ws.Range("H132").Value = 4113.476
$ws.Range("I132").Value = 3470.2856
$ws.Range("J132").Value = 5399.857
$ws.Range("K132").Value = 10410.8568
$ws.Range("L132").Value = 16199.571
$ws.Range("M132").Value = -7880.856800000001
$ws.Range("N132").Value = -21259.571


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3560.35
$ws.Range("I7").Value = 2088.6667
$ws.Range("K7").Value = 2088.6667
$ws.Range("M7").Value = -1976.6667

$ws.Range("H22").Value = 691.8333
$ws.Range("I22").Value = 457.42856
$ws.Range("K22").Value = 457.42856
$ws.Range("M22").Value = -162.42856

$ws.Range("H27").Value = 691.8333
$ws.Range("I27").Value = 457.42856
$ws.Range("K27").Value = 457.42856
$ws.Range("M27").Value = -350.42856

$ws.Range("H68").Value = 1720.9474
$ws.Range("I68").Value = 1833.5454
$ws.Range("J68").Value = 1566.125
$ws.Range("K68").Value = 1833.5454
$ws.Range("L68").Value = 1566.125
$ws.Range("M68").Value = -1084.5454
$ws.Range("N68").Value = -3064.125

$ws.Range("H71").Value = 1720.9474
$ws.Range("I71").Value = 1833.5454
$ws.Range("J71").Value = 1566.125
$ws.Range("K71").Value = 9167.726999999999
$ws.Range("L71").Value = 7830.625
$ws.Range("M71").Value = -5423.726999999999
$ws.Range("N71").Value = -15318.625

$ws.Range("H122").Value = 5655.524
$ws.Range("I122").Value = 4032.6667
$ws.Range("K122").Value = 12098.0001
$ws.Range("M122").Value = -9648.000100000001

$ws.Range("H126").Value = 3560.35
$ws.Range("I126").Value = 2088.6667
$ws.Range("K126").Value = 6266.000100000001
$ws.Range("M126").Value = -3796.000100000001


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 49996.5
$ws.Range("J109").Value = 49996.5
$ws.Range("L109").Value = 49996.5
$ws.Range("N109").Value = -52770.5

$ws.Range("H126").Value = 1789.9615
$ws.Range("I126").Value = 1867.25
$ws.Range("K126").Value = 5601.75
$ws.Range("M126").Value = -3131.75

$ws.Range("H136").Value = 2319.158
$ws.Range("I136").Value = 1178.4572
$ws.Range("J136").Value = 4133.909
$ws.Range("K136").Value = 3535.3716
$ws.Range("L136").Value = 12401.727
$ws.Range("M136").Value = -985.3716000000004
$ws.Range("N136").Value = -17501.727

